# "arrow + bugfix plank" --------------------------------------------------
# Recomputes the "plank" cut list (rows 7-13 in the original sheet): several
# lengths/widths/counts change, one extra pair of plank rows is added, and
# the row block that used to span rows 7-13 (plus the "scharnier"/"slot"
# rows that follow it) grows to rows 7-15 / 16-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10:C13 is about to be split into two separate merges (C10:C11, C12:C15),
# so break it apart before touching any values underneath it.
$ws.Range("C10:C13").UnMerge()

# --- Updated values for the existing "balk" rows (2-5) -------------------
$ws.Range("D2").Value = 69.59999999999999
$ws.Range("E2").Value = 6

$ws.Range("D3").Value = 75.59999999999999

$ws.Range("D4").Value = 93.3
$ws.Range("E4").Value = 4

$ws.Range("D5").Value = 184.6

# --- Updated values for the existing "plank" rows (7-13) -----------------
$ws.Range("D7").Value = 195

$ws.Range("D8").Value = 195

$ws.Range("C9").Value = 14.8
$ws.Range("D9").Value = 194.6
$ws.Range("E9").Value = 4

$ws.Range("C10").Value = 17.8
$ws.Range("D10").Value = 195
$ws.Range("E10").Value = 4

$ws.Range("D11").Value = 195.6
$ws.Range("E11").Value = 4

$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 79.40000000000001
$ws.Range("E12").Value = 6

$ws.Range("D13").Value = 194.6
$ws.Range("E13").Value = 4

# --- Insert two new "plank" rows, pushing "scharnier"/"slot" down --------
$ws.Rows("14:15").Insert()

# Format the new rows like the rest of the bordered/centered block (style
# index 1 in the original workbook: thin box border, bold, center/top).
$ws.Range("A14:D15").Font.Bold = $true
$ws.Range("A14:D15").Borders.LineStyle = 1
$ws.Range("A14:D15").HorizontalAlignment = -4108
$ws.Range("A14:D15").VerticalAlignment = -4160

$ws.Range("D14").Value = 195
$ws.Range("E14").Value = 20

$ws.Range("D15").Value = 195.6
$ws.Range("E15").Value = 4

# --- Extend the merges that used to stop at row 13 to the new row 15 -----
$ws.Range("A7:A15").Merge()
$ws.Range("B7:B15").Merge()
$ws.Range("C10:C11").Merge()
$ws.Range("C12:C15").Merge()

# Merging repaints the merged block's outer border only; restore the
# uniform per-cell thin border / bold / center-top formatting so every
# cell in the block keeps looking like the rest of the table.
$ws.Range("A7:D15").Borders.LineStyle = 1
$ws.Range("A7:D15").Font.Bold = $true
$ws.Range("A7:D15").HorizontalAlignment = -4108
$ws.Range("A7:D15").VerticalAlignment = -4160
